$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.412.74"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.246.43"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.07"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.60"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.95"
$ws.Range("E10").Value = "  +9.05%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.581.84"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.56"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.857"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.238.07"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.217.29"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.15"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.24"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.11"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.15"
$ws.Range("E24").Value = "  +28.93%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("E26").Value = "  +2.77%  "
$ws.Range("E27").Value = "  -2.86%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.52"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.90"
$ws.Range("E34").Value = "  -5.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.35"
$ws.Range("E35").Value = "  +11.08%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0315"
$ws.Range("E38").Value = "  +6.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.94"
$ws.Range("E39").Value = "  +5.00%  "
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.79"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.08"
$ws.Range("E42").Value = "  +6.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.202"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.73"
$ws.Range("E44").Value = "  -7.90%  "
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("E51").Value = "  +0.96%  "
